{"js": "// Replace each arithmetic-problem cell's text with its new value.\n// `pairs[i] = [oldText, newText]` corresponds to the table cell at\n// row = floor(i / 5), col = i % 5 (the table is 20 rows x 5 cols,\n// filled left-to-right, top-to-bottom -- matching the order the\n// problems appear in the document).\nconst pairs = [\n  [\"24-12=\", \"90-65=\"],\n  [\"39+23=\", \"31+6=\"],\n  [\"37-28=\", \"17+44=\"],\n  [\"5+59=\", \"61+33=\"],\n  [\"87-26=\", \"41-4=\"],\n  [\"74-3=\", \"83-71=\"],\n  [\"13+81=\", \"82+0=\"],\n  [\"18-6=\", \"27+46=\"],\n  [\"34-5=\", \"85-36=\"],\n  [\"8+43=\", \"68-44=\"],\n  [\"35+57=\", \"27-7=\"],\n  [\"79-55=\", \"35+59=\"],\n  [\"75-29=\", \"43-12=\"],\n  [\"59-15=\", \"25+51=\"],\n  [\"27+60=\", \"88-74=\"],\n  [\"34+43=\", \"91+0=\"],\n  [\"58+35=\", \"41+13=\"],\n  [\"42+15=\", \"82-59=\"],\n  [\"86-16=\", \"92-87=\"],\n  [\"64-44=\", \"89-72=\"],\n  [\"38+40=\", \"75+12=\"],\n  [\"75-0=\", \"39-4=\"],\n  [\"59+5=\", \"84-49=\"],\n  [\"20+41=\", \"98-38=\"],\n  [\"7+73=\", \"32-13=\"],\n  [\"64-63=\", \"62-31=\"],\n  [\"79-38=\", \"20-15=\"],\n  [\"92-67=\", \"35-1=\"],\n  [\"65+22=\", \"13-7=\"],\n  [\"69-44=\", \"71-38=\"],\n  [\"50-5=\", \"72-59=\"],\n  [\"16-11=\", \"1+0=\"],\n  [\"70-61=\", \"48-19=\"],\n  [\"45-10=\", \"43+31=\"],\n  [\"66-34=\", \"86-46=\"],\n  [\"52-0=\", \"25+46=\"],\n  [\"13+37=\", \"46-4=\"],\n  [\"28+52=\", \"76-10=\"],\n  [\"58+3=\", \"93-87=\"],\n  [\"14-9=\", \"68-39=\"],\n  [\"39+21=\", \"25+18=\"],\n  [\"13+30=\", \"83-7=\"],\n  [\"99-47=\", \"89-76=\"],\n  [\"11+50=\", \"72-48=\"],\n  [\"84-22=\", \"0+51=\"],\n  [\"79-29=\", \"22+19=\"],\n  [\"53+14=\", \"40-36=\"],\n  [\"71-65=\", \"51-29=\"],\n  [\"74+20=\", \"52+11=\"],\n  [\"69+4=\", \"11+16=\"],\n  [\"13+73=\", \"26+27=\"],\n  [\"87-44=\", \"62+26=\"],\n  [\"95-24=\", \"56+26=\"],\n  [\"65-30=\", \"58-1=\"],\n  [\"58-8=\", \"30+21=\"],\n  [\"16+80=\", \"16+53=\"],\n  [\"81-31=\", \"32+31=\"],\n  [\"50-6=\", \"73-30=\"],\n  [\"31-14=\", \"22+32=\"],\n  [\"53+22=\", \"70+12=\"],\n  [\"43-16=\", \"22+24=\"],\n  [\"38+0=\", \"18+51=\"],\n  [\"48-43=\", \"91-26=\"],\n  [\"9+16=\", \"80+3=\"],\n  [\"49-36=\", \"63+5=\"],\n  [\"67-41=\", \"24+32=\"],\n  [\"81-26=\", \"23+72=\"],\n  [\"64+3=\", \"87-47=\"],\n  [\"63+6=\", \"53+33=\"],\n  [\"36+35=\", \"47+8=\"],\n  [\"7+17=\", \"96-90=\"],\n  [\"56+21=\", \"55-37=\"],\n  [\"9+67=\", \"57-12=\"],\n  [\"75+16=\", \"85-61=\"],\n  [\"94-73=\", \"11+84=\"],\n  [\"14-3=\", \"80-1=\"],\n  [\"13+23=\", \"6+1=\"],\n  [\"65+17=\", \"51-27=\"],\n  [\"90-72=\", \"83-59=\"],\n  [\"63+26=\", \"75-47=\"],\n  [\"37+32=\", \"15+84=\"],\n  [\"28+15=\", \"65-57=\"],\n  [\"77-13=\", \"41+19=\"],\n  [\"11+59=\", \"48-19=\"],\n  [\"47-4=\", \"62-25=\"],\n  [\"75-47=\", \"13+86=\"],\n  [\"77-51=\", \"78-35=\"],\n  [\"76+14=\", \"38-11=\"],\n  [\"41-35=\", \"94-78=\"],\n  [\"45+30=\", \"98-22=\"],\n  [\"21+66=\", \"84-11=\"],\n  [\"21+58=\", \"1+41=\"],\n  [\"82-48=\", \"99-73=\"],\n  [\"11+60=\", \"31-22=\"],\n  [\"76+3=\", \"70-20=\"],\n  [\"66-21=\", \"18+54=\"],\n  [\"33+36=\", \"87+2=\"],\n  [\"66+24=\", \"44-26=\"],\n  [\"45+9=\", \"84-38=\"],\n  [\"57-39=\", \"19+39=\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst numCols = 5;\nif (pairs.length !== table.rowCount * numCols) {\n  throw new Error(\n    \"Unexpected cell count: \" + table.rowCount * numCols + \" vs \" + pairs.length\n  );\n}\n\n// Search each cell individually (not the whole document) for its old\n// text, then replace it. Scoping the search to the owning cell keeps\n// the operation correct even though some \"new\" values equal other\n// cells' \"old\" values elsewhere in the table.\nconst searchResults = [];\nfor (let i = 0; i < pairs.length; i++) {\n  const row = Math.floor(i / numCols);\n  const col = i % numCols;\n  const [oldText] = pairs[i];\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const results = searchResults[i];\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + oldText + \"' in cell \" + i +\n        \", got \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace each arithmetic-problem cell's text with its new value.\n# $pairs[$i] = @(oldText, newText) corresponds to the table cell at\n# row = floor(i / 5) + 1, col = (i % 5) + 1 (the table is 20 rows x 5\n# cols, 1-indexed, filled left-to-right top-to-bottom -- matching the\n# order the problems appear in the document).\n$pairs = @(\n    @(\"24-12=\", \"90-65=\"),\n    @(\"39+23=\", \"31+6=\"),\n    @(\"37-28=\", \"17+44=\"),\n    @(\"5+59=\", \"61+33=\"),\n    @(\"87-26=\", \"41-4=\"),\n    @(\"74-3=\", \"83-71=\"),\n    @(\"13+81=\", \"82+0=\"),\n    @(\"18-6=\", \"27+46=\"),\n    @(\"34-5=\", \"85-36=\"),\n    @(\"8+43=\", \"68-44=\"),\n    @(\"35+57=\", \"27-7=\"),\n    @(\"79-55=\", \"35+59=\"),\n    @(\"75-29=\", \"43-12=\"),\n    @(\"59-15=\", \"25+51=\"),\n    @(\"27+60=\", \"88-74=\"),\n    @(\"34+43=\", \"91+0=\"),\n    @(\"58+35=\", \"41+13=\"),\n    @(\"42+15=\", \"82-59=\"),\n    @(\"86-16=\", \"92-87=\"),\n    @(\"64-44=\", \"89-72=\"),\n    @(\"38+40=\", \"75+12=\"),\n    @(\"75-0=\", \"39-4=\"),\n    @(\"59+5=\", \"84-49=\"),\n    @(\"20+41=\", \"98-38=\"),\n    @(\"7+73=\", \"32-13=\"),\n    @(\"64-63=\", \"62-31=\"),\n    @(\"79-38=\", \"20-15=\"),\n    @(\"92-67=\", \"35-1=\"),\n    @(\"65+22=\", \"13-7=\"),\n    @(\"69-44=\", \"71-38=\"),\n    @(\"50-5=\", \"72-59=\"),\n    @(\"16-11=\", \"1+0=\"),\n    @(\"70-61=\", \"48-19=\"),\n    @(\"45-10=\", \"43+31=\"),\n    @(\"66-34=\", \"86-46=\"),\n    @(\"52-0=\", \"25+46=\"),\n    @(\"13+37=\", \"46-4=\"),\n    @(\"28+52=\", \"76-10=\"),\n    @(\"58+3=\", \"93-87=\"),\n    @(\"14-9=\", \"68-39=\"),\n    @(\"39+21=\", \"25+18=\"),\n    @(\"13+30=\", \"83-7=\"),\n    @(\"99-47=\", \"89-76=\"),\n    @(\"11+50=\", \"72-48=\"),\n    @(\"84-22=\", \"0+51=\"),\n    @(\"79-29=\", \"22+19=\"),\n    @(\"53+14=\", \"40-36=\"),\n    @(\"71-65=\", \"51-29=\"),\n    @(\"74+20=\", \"52+11=\"),\n    @(\"69+4=\", \"11+16=\"),\n    @(\"13+73=\", \"26+27=\"),\n    @(\"87-44=\", \"62+26=\"),\n    @(\"95-24=\", \"56+26=\"),\n    @(\"65-30=\", \"58-1=\"),\n    @(\"58-8=\", \"30+21=\"),\n    @(\"16+80=\", \"16+53=\"),\n    @(\"81-31=\", \"32+31=\"),\n    @(\"50-6=\", \"73-30=\"),\n    @(\"31-14=\", \"22+32=\"),\n    @(\"53+22=\", \"70+12=\"),\n    @(\"43-16=\", \"22+24=\"),\n    @(\"38+0=\", \"18+51=\"),\n    @(\"48-43=\", \"91-26=\"),\n    @(\"9+16=\", \"80+3=\"),\n    @(\"49-36=\", \"63+5=\"),\n    @(\"67-41=\", \"24+32=\"),\n    @(\"81-26=\", \"23+72=\"),\n    @(\"64+3=\", \"87-47=\"),\n    @(\"63+6=\", \"53+33=\"),\n    @(\"36+35=\", \"47+8=\"),\n    @(\"7+17=\", \"96-90=\"),\n    @(\"56+21=\", \"55-37=\"),\n    @(\"9+67=\", \"57-12=\"),\n    @(\"75+16=\", \"85-61=\"),\n    @(\"94-73=\", \"11+84=\"),\n    @(\"14-3=\", \"80-1=\"),\n    @(\"13+23=\", \"6+1=\"),\n    @(\"65+17=\", \"51-27=\"),\n    @(\"90-72=\", \"83-59=\"),\n    @(\"63+26=\", \"75-47=\"),\n    @(\"37+32=\", \"15+84=\"),\n    @(\"28+15=\", \"65-57=\"),\n    @(\"77-13=\", \"41+19=\"),\n    @(\"11+59=\", \"48-19=\"),\n    @(\"47-4=\", \"62-25=\"),\n    @(\"75-47=\", \"13+86=\"),\n    @(\"77-51=\", \"78-35=\"),\n    @(\"76+14=\", \"38-11=\"),\n    @(\"41-35=\", \"94-78=\"),\n    @(\"45+30=\", \"98-22=\"),\n    @(\"21+66=\", \"84-11=\"),\n    @(\"21+58=\", \"1+41=\"),\n    @(\"82-48=\", \"99-73=\"),\n    @(\"11+60=\", \"31-22=\"),\n    @(\"76+3=\", \"70-20=\"),\n    @(\"66-21=\", \"18+54=\"),\n    @(\"33+36=\", \"87+2=\"),\n    @(\"66+24=\", \"44-26=\"),\n    @(\"45+9=\", \"84-38=\"),\n    @(\"57-39=\", \"19+39=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$numCols = 5\n\nif ($pairs.Count -ne ($t.Rows.Count * $numCols)) {\n    throw \"Unexpected cell count: expected $($t.Rows.Count * $numCols), got $($pairs.Count)\"\n}\n\n# Address each cell by its fixed (row, col) position rather than by\n# searching document-wide for the old text: several old/new values\n# coincide across different cells, so a global Find/Replace could\n# touch the wrong cell. Excluding the trailing cell-mark character\n# (End - 1) keeps the run's formatting (font/size) intact.\nfor ($i = 0; $i -lt $pairs.Count; $i++) {\n    $row = [int][Math]::Floor($i / $numCols) + 1\n    $col = ($i % $numCols) + 1\n    $oldText = $pairs[$i][0]\n    $newText = $pairs[$i][1]\n\n    $cell = $t.Cell($row, $col)\n    $cellStart = $cell.Range.Start\n    $cellEnd = $cell.Range.End\n    $textRng = $d.Range($cellStart, $cellEnd - 1)\n\n    if ($textRng.Text -ne $oldText) {\n        throw \"Cell ($row,$col) text mismatch: expected '$oldText', found '$($textRng.Text)'\"\n    }\n\n    $textRng.Text = $newText\n}\n\nWrite-Output \"done\"\n"}
